$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (B, C, D) keep their values as plain text
# rather than being auto-converted to numbers by Excel smart-parsing.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.282.54'
$ws.Range("E2").Value = '  +1.86%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.914.41'

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.56'
$ws.Range("E5").Value = '  +0.98%  '

$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4629'
$ws.Range("E7").Value = '  +0.79%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3947'
$ws.Range("E8").Value = '  +2.27%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.75'
$ws.Range("E9").Value = '  +1.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07944'
$ws.Range("E10").Value = '  +0.99%  '

$ws.Range("E11").Value = '  +0.74%  '

$ws.Range("E12").Value = '  +2.97%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.917.03'
$ws.Range("E13").Value = '  +0.40%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.114'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.773'
$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06941'
$ws.Range("E16").Value = '  -0.52%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.63'
$ws.Range("E17").Value = '  +0.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.003'
$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001006'
$ws.Range("E19").Value = '  +0.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.12'
$ws.Range("E20").Value = '  +1.71%  '

$ws.Range("E21").Value = '  -0.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29.302.24'
$ws.Range("E22").Value = '  +1.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.360'
$ws.Range("E23").Value = '  +1.54%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  +0.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.155.14'
$ws.Range("E25").Value = '  +1.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.055'
$ws.Range("E26").Value = '  -3.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.04'
$ws.Range("E27").Value = '  +2.36%  '

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.204'
$ws.Range("E28").Value = '  +7.51%  '

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.48'
$ws.Range("E29").Value = '  +1.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.996'
$ws.Range("E30").Value = '  +2.33%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '118.92'
$ws.Range("E31").Value = '  -0.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09395'
$ws.Range("E32").Value = '  +0.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9261'
$ws.Range("E33").Value = '  +0.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.362'
$ws.Range("E34").Value = '  +1.23%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.358'
$ws.Range("E35").Value = '  +1.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.275'
$ws.Range("E36").Value = '  -1.02%  '

$ws.Range("E37").Value = '  +4.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05842'
$ws.Range("E38").Value = '  +1.58%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02111'
$ws.Range("E39").Value = '  +1.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.979'
$ws.Range("E40").Value = '  +3.41%  '

$ws.Range("E41").Value = '  -0.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5755'
$ws.Range("E42").Value = '  +2.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1803'
$ws.Range("E43").Value = '  +0.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.960'
$ws.Range("E44").Value = '  +0.94%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.315'
$ws.Range("E45").Value = '  +8.33%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.02'
$ws.Range("E46").Value = '  +1.94%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5427'
$ws.Range("E47").Value = '  +2.84%  '

$ws.Range("E48").Value = '  -1.62%  '

$ws.Range("E49").Value = '  +3.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.563'
$ws.Range("E50").Value = '  +6.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '113.56'
$ws.Range("E51").Value = '  +0.06%  '
